$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(112, 8).Value = 5315.222
$ws.Cells.Item(112, 10).Value = 5539.647
$ws.Cells.Item(112, 12).Value = 16618.941
$ws.Cells.Item(112, 14).Value = -18834.941

$ws.Cells.Item(113, 8).Value = 6164.75
$ws.Cells.Item(113, 9).Value = 8124.75
$ws.Cells.Item(113, 10).Value = 5184.75
$ws.Cells.Item(113, 11).Value = 8124.75
$ws.Cells.Item(113, 12).Value = 5184.75
$ws.Cells.Item(113, 13).Value = -4870.75
$ws.Cells.Item(113, 14).Value = -11692.75

$ws.Cells.Item(129, 8).Value = 994.0833
$ws.Cells.Item(129, 9).Value = 994.0833
$ws.Cells.Item(129, 11).Value = 2982.2499
$ws.Cells.Item(129, 13).Value = 2017.7501

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 4521.8887
$ws.Cells.Item(2, 9).Value = 3499
$ws.Cells.Item(2, 11).Value = 3499
$ws.Cells.Item(2, 13).Value = -3386

$ws.Cells.Item(5, 8).Value = 2630.0952
$ws.Cells.Item(5, 9).Value = 2281.2144
$ws.Cells.Item(5, 10).Value = 3327.8572
$ws.Cells.Item(5, 11).Value = 2281.2144
$ws.Cells.Item(5, 12).Value = 3327.8572
$ws.Cells.Item(5, 13).Value = -2169.2144
$ws.Cells.Item(5, 14).Value = -3551.8572

$ws.Cells.Item(32, 8).Value = 7258.357
$ws.Cells.Item(32, 9).Value = 733.5319
$ws.Cells.Item(32, 10).Value = 41332.445
$ws.Cells.Item(32, 11).Value = 733.5319
$ws.Cells.Item(32, 12).Value = 41332.445
$ws.Cells.Item(32, 13).Value = -446.5319
$ws.Cells.Item(32, 14).Value = -41906.445

$ws.Cells.Item(74, 8).Value = 2074.4167
$ws.Cells.Item(74, 9).Value = 2118.9
$ws.Cells.Item(74, 11).Value = 2118.9
$ws.Cells.Item(74, 13).Value = -1244.9

$ws.Cells.Item(77, 8).Value = 2074.4167
$ws.Cells.Item(77, 9).Value = 2118.9
$ws.Cells.Item(77, 11).Value = 10594.5
$ws.Cells.Item(77, 13).Value = -6226.5

$ws.Cells.Item(116, 8).Value = 4521.8887
$ws.Cells.Item(116, 9).Value = 3499
$ws.Cells.Item(116, 11).Value = 3499
$ws.Cells.Item(116, 13).Value = -1205

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 4521.8887
$ws.Cells.Item(3, 9).Value = 3499
$ws.Cells.Item(3, 11).Value = 3499
$ws.Cells.Item(3, 13).Value = -3385

$ws.Cells.Item(4, 8).Value = 2630.0952
$ws.Cells.Item(4, 9).Value = 2281.2144
$ws.Cells.Item(4, 10).Value = 3327.8572
$ws.Cells.Item(4, 11).Value = 2281.2144
$ws.Cells.Item(4, 12).Value = 3327.8572
$ws.Cells.Item(4, 13).Value = -2166.2144
$ws.Cells.Item(4, 14).Value = -3557.8572

$ws.Cells.Item(80, 8).Value = 1618.8235
$ws.Cells.Item(80, 9).Value = 1640.5555
$ws.Cells.Item(80, 10).Value = 1594.375
$ws.Cells.Item(80, 11).Value = 1640.5555
$ws.Cells.Item(80, 12).Value = 1594.375
$ws.Cells.Item(80, 13).Value = -642.5554999999999
$ws.Cells.Item(80, 14).Value = -3590.375

$ws.Cells.Item(83, 8).Value = 1618.8235
$ws.Cells.Item(83, 9).Value = 1640.5555
$ws.Cells.Item(83, 10).Value = 1594.375
$ws.Cells.Item(83, 11).Value = 8202.7775
$ws.Cells.Item(83, 12).Value = 7971.875
$ws.Cells.Item(83, 13).Value = -3210.7775
$ws.Cells.Item(83, 14).Value = -17955.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 12731.308
$ws.Cells.Item(16, 9).Value = 1115.8572
$ws.Cells.Item(16, 10).Value = 26282.666
$ws.Cells.Item(16, 11).Value = 1115.8572
$ws.Cells.Item(16, 12).Value = 26282.666
$ws.Cells.Item(16, 13).Value = -828.8571999999999
$ws.Cells.Item(16, 14).Value = -26856.666

$ws.Cells.Item(62, 8).Value = 20000
$ws.Cells.Item(62, 9).Value = 20000
$ws.Cells.Item(62, 11).Value = 20000
$ws.Cells.Item(62, 13).Value = -19376

$ws.Cells.Item(65, 8).Value = 20000
$ws.Cells.Item(65, 9).Value = 20000
$ws.Cells.Item(65, 11).Value = 100000
$ws.Cells.Item(65, 13).Value = -96880

$ws.Cells.Item(105, 8).Value = 1145
$ws.Cells.Item(105, 9).Value = 971.1111
$ws.Cells.Item(105, 11).Value = 971.1111
$ws.Cells.Item(105, 13).Value = 775.8889

$ws.Cells.Item(107, 8).Value = 1668.5
$ws.Cells.Item(107, 9).Value = 2011
$ws.Cells.Item(107, 11).Value = 2011
$ws.Cells.Item(107, 13).Value = -91

$ws.Cells.Item(113, 8).Value = 12731.308
$ws.Cells.Item(113, 9).Value = 1115.8572
$ws.Cells.Item(113, 10).Value = 26282.666
$ws.Cells.Item(113, 11).Value = 1115.8572
$ws.Cells.Item(113, 12).Value = 26282.666
$ws.Cells.Item(113, 13).Value = 1054.1428
$ws.Cells.Item(113, 14).Value = -30622.666

$ws.Cells.Item(134, 8).Value = 4429.0347
$ws.Cells.Item(134, 9).Value = 4386.7407
$ws.Cells.Item(134, 11).Value = 13160.2221
$ws.Cells.Item(134, 13).Value = -10625.2221

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 1087.5883
$ws.Cells.Item(2, 10).Value = 582.75
$ws.Cells.Item(2, 12).Value = 3496.5
$ws.Cells.Item(2, 14).Value = -3722.5

$ws.Cells.Item(139, 8).Value = 6911.6523
$ws.Cells.Item(139, 9).Value = 3092.4546
$ws.Cells.Item(139, 11).Value = 9277.363799999999
$ws.Cells.Item(139, 13).Value = -4137.363799999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(108, 8).Value = 70310
$ws.Cells.Item(108, 9).Value = 40621
$ws.Cells.Item(108, 11).Value = 40621
$ws.Cells.Item(108, 13).Value = -36781

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 2629.04
$ws.Cells.Item(22, 9).Value = 2443.5386
$ws.Cells.Item(22, 11).Value = 2443.5386
$ws.Cells.Item(22, 13).Value = -2148.5386

$ws.Cells.Item(27, 8).Value = 2629.04
$ws.Cells.Item(27, 9).Value = 2443.5386
$ws.Cells.Item(27, 11).Value = 2443.5386
$ws.Cells.Item(27, 13).Value = -2336.5386

$ws.Cells.Item(40, 8).Value = 10507.5
$ws.Cells.Item(40, 9).Value = 7720.5557
$ws.Cells.Item(40, 10).Value = 18868.334
$ws.Cells.Item(40, 11).Value = 7720.5557
$ws.Cells.Item(40, 12).Value = 18868.334
$ws.Cells.Item(40, 13).Value = -7584.5557
$ws.Cells.Item(40, 14).Value = -19140.334

$ws.Cells.Item(46, 8).Value = 5758.1763
$ws.Cells.Item(46, 9).Value = 7099.857
$ws.Cells.Item(46, 10).Value = 4819
$ws.Cells.Item(46, 11).Value = 7099.857
$ws.Cells.Item(46, 12).Value = 4819
$ws.Cells.Item(46, 13).Value = -6911.857
$ws.Cells.Item(46, 14).Value = -5195

$ws.Cells.Item(68, 8).Value = 2697.077
$ws.Cells.Item(68, 9).Value = 1506.875
$ws.Cells.Item(68, 10).Value = 4601.4
$ws.Cells.Item(68, 11).Value = 1506.875
$ws.Cells.Item(68, 12).Value = 4601.4
$ws.Cells.Item(68, 13).Value = -757.875
$ws.Cells.Item(68, 14).Value = -6099.4

$ws.Cells.Item(71, 8).Value = 2697.077
$ws.Cells.Item(71, 9).Value = 1506.875
$ws.Cells.Item(71, 10).Value = 4601.4
$ws.Cells.Item(71, 11).Value = 7534.375
$ws.Cells.Item(71, 12).Value = 23007
$ws.Cells.Item(71, 13).Value = -3790.375
$ws.Cells.Item(71, 14).Value = -30495

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 4206.1875
$ws.Cells.Item(132, 9).Value = 6043
$ws.Cells.Item(132, 10).Value = 2777.5557
$ws.Cells.Item(132, 11).Value = 18129
$ws.Cells.Item(132, 12).Value = 8332.667099999999
$ws.Cells.Item(132, 13).Value = -15599
$ws.Cells.Item(132, 14).Value = -13392.6671
